$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix cell styles/types for cells that switch between text-placeholder and numeric ---
$ws.Range("F15").Copy($ws.Range("C15"))
$ws.Range("F15").Copy($ws.Range("D15"))
$ws.Range("H15").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("C14").Copy($ws.Range("E22"))
$ws.Range("F15").Copy($ws.Range("C27"))
$ws.Range("F15").Copy($ws.Range("D27"))
$ws.Range("H15").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("C14").Copy($ws.Range("G33"))
$ws.Range("C14").Copy($ws.Range("H33"))
$ws.Range("H15").Copy($ws.Range("L33"))

# --- Step 2: set final cell values ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"
$ws.Range("M14").Value = 25
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = -27.272727272727
$ws.Range("L15").Value = 45.454545454545
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -30.434782608695
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 56
$ws.Range("G16").Value = 57
$ws.Range("H16").Value = -1.754385964912
$ws.Range("I16").Value = 236
$ws.Range("J16").Value = 215
$ws.Range("K16").Value = 9.767441860465
$ws.Range("L16").Value = 22.916666666666
$ws.Range("M16").Value = 37.209302325581
$ws.Range("N16").Value = -67.891156462585
$ws.Range("C17").Value = 27
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 90
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 390
$ws.Range("J17").Value = 373
$ws.Range("K17").Value = 4.557640750670
$ws.Range("L17").Value = 17.469879518072
$ws.Range("M17").Value = 124.137931034483
$ws.Range("N17").Value = -4.645476772616
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 166.666666666667
$ws.Range("F18").Value = 32
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = -8.571428571428
$ws.Range("I18").Value = 143
$ws.Range("J18").Value = 129
$ws.Range("K18").Value = 10.852713178294
$ws.Range("L18").Value = -1.379310344827
$ws.Range("M18").Value = 70.238095238095
$ws.Range("N18").Value = -73.809523809523
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 48.936170212766
$ws.Range("I19").Value = 367
$ws.Range("J19").Value = 270
$ws.Range("K19").Value = 35.925925925925
$ws.Range("L19").Value = 31.541218637992
$ws.Range("M19").Value = 144.666666666667
$ws.Range("N19").Value = 37.969924812030
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -59.259259259259
$ws.Range("I20").Value = 80
$ws.Range("J20").Value = 144
$ws.Range("K20").Value = -44.444444444444
$ws.Range("L20").Value = -38.931297709923
$ws.Range("M20").Value = 86.046511627907
$ws.Range("N20").Value = -69.111969111969
$ws.Range("C21").Value = 69
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = 30.188679245283
$ws.Range("F21").Value = 251
$ws.Range("G21").Value = 261
$ws.Range("H21").Value = -3.831417624521
$ws.Range("I21").Value = 1237
$ws.Range("J21").Value = 1158
$ws.Range("K21").Value = 6.822107081174
$ws.Range("L21").Value = 12.762078395624
$ws.Range("M21").Value = 94.803149606299
$ws.Range("N21").Value = -45.458553791887
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 31
$ws.Range("K22").Value = 82.352941176470
$ws.Range("L22").Value = -34.042553191489
$ws.Range("M22").Value = 29.166666666666
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 55
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 195
$ws.Range("J23").Value = 208
$ws.Range("K23").Value = -6.25
$ws.Range("L23").Value = 33.561643835616
$ws.Range("M23").Value = 77.272727272727
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -14.814814814814
$ws.Range("F24").Value = 112
$ws.Range("H24").Value = -2.608695652173
$ws.Range("I24").Value = 635
$ws.Range("J24").Value = 631
$ws.Range("K24").Value = 0.633914421553
$ws.Range("L24").Value = 5.481727574750
$ws.Range("M24").Value = 28.803245436105
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 54.166666666666
$ws.Range("I25").Value = 254
$ws.Range("J25").Value = 233
$ws.Range("K25").Value = 9.012875536480
$ws.Range("L25").Value = -4.868913857677
$ws.Range("C26").Value = 29
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 70.588235294117
$ws.Range("G26").Value = 93
$ws.Range("H26").Value = 8.602150537634
$ws.Range("I26").Value = 484
$ws.Range("J26").Value = 446
$ws.Range("K26").Value = 8.520179372197
$ws.Range("L26").Value = 26.041666666666
$ws.Range("M26").Value = 14.964370546318
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = -28.125
$ws.Range("L27").Value = 35.294117647058
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 185.714285714286
$ws.Range("I28").Value = 66
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 26.923076923076
$ws.Range("L28").Value = 144.444444444444
$ws.Range("C29").Value = "0"
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -41.666666666666
$ws.Range("M29").Value = -41.666666666666
$ws.Range("N29").Value = -80.821917808219
$ws.Range("C30").Value = "0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("L30").Value = -27.777777777777
$ws.Range("M30").Value = -23.529411764705
$ws.Range("N30").Value = -79.365079365079
$ws.Range("G33").Value = "0"
$ws.Range("H33").Value = "***.*"
$ws.Range("L33").Value = 0

# --- Step 3: column E width (bestFit shrink due to shorter % values) ---
$ws.Columns("E").ColumnWidth = 6.168446
